$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A72").Value = "27-11-2025"
$ws.Range("B72").Value = "The price of gold in India today is ₹12,775 per gram for 24 karat gold, ₹11,710 per gram for 22 karat gold and ₹9,581 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A72").Borders.LineStyle = 1
$ws.Range("B72").Borders.LineStyle = 1
$ws.Range("B72").WrapText = $true
